# Fix Ubuntu 14 AMI ids for AWS regions (row 4 on Sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Existing AMI ids (D4, F4, H4, I4) get new values, and two new region
# columns (G4 = us-west-2, J4 = eu-west-2) are populated for the first time.
$ws.Range("D4").Value = "ami-9dde7f8b"
$ws.Range("F4").Value = "ami-9d772efd"
$ws.Range("G4").Value = "ami-0e2aa66e"
$ws.Range("H4").Value = "ami-115d7777"
$ws.Range("I4").Value = "ami-6039ed0f"
$ws.Range("J4").Value = "ami-c29184a6"

# Move the cursor/selection from J10 to J9 (matches the saved view state).
$null = $ws.Range("J9").Select()
